$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the tag value cells per the diff
$ws.Range("D2").Value = "Value08"
$ws.Range("D3").Value = "Value09"
$ws.Range("A2").Value = "Default - Microsoft Azure Sponsorship 3"

# Update the selected cell in the sheet view
$ws.Range("A3").Select()
